$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.87"
$ws.Range("E2").Value = "'2.32%"
$ws.Range("D3").Value = "'31.98"
$ws.Range("E3").Value = "'2.68%"
$ws.Range("D4").Value = "'5.025"
$ws.Range("E4").Value = "'1.90%"
$ws.Range("D5").Value = "'0.07782"
$ws.Range("E5").Value = "'6.05%"
$ws.Range("D6").Value = "'2.301"
$ws.Range("E6").Value = "'0.02%"
$ws.Range("D7").Value = "'7.937"
$ws.Range("E7").Value = "'3.15%"
$ws.Range("D8").Value = "'0.9320"
$ws.Range("E8").Value = "'2.15%"
$ws.Range("D9").Value = "'0.1016"
$ws.Range("E9").Value = "'24.62%"
$ws.Range("D10").Value = "'0.1761"
$ws.Range("E10").Value = "'4.45%"
$ws.Range("D11").Value = "'0.08421"
$ws.Range("E11").Value = "'1.63%"
$ws.Range("D12").Value = "'0.03312"
$ws.Range("E12").Value = "'6.68%"
$ws.Range("D13").Value = "'0.09902"
$ws.Range("E13").Value = "'-1.54%"
$ws.Range("D14").Value = "'0.001472"
$ws.Range("E14").Value = "'-2.89%"
$ws.Range("D15").Value = "'0.005702"
$ws.Range("E15").Value = "'-0.22%"
$ws.Range("D16").Value = "'3.495"
$ws.Range("D17").Value = "'3.850"
$ws.Range("E17").Value = "'2.94%"
$ws.Range("D18").Value = "'2.191"
$ws.Range("E18").Value = "'5.45%"
$ws.Range("D19").Value = "'0.3359"
$ws.Range("E19").Value = "'0.90%"
$ws.Range("D20").Value = "'0.1341"
$ws.Range("E20").Value = "'2.87%"
$ws.Range("D21").Value = "'4.293"
$ws.Range("E21").Value = "'8.05%"
$ws.Range("D22").Value = "'0.2081"
$ws.Range("E22").Value = "'-1.00%"
$ws.Range("D23").Value = "'0.04614"
$ws.Range("E23").Value = "'1.23%"
$ws.Range("D24").Value = "'0.001214"
$ws.Range("E24").Value = "'0.28%"
$ws.Range("D25").Value = "'0.004389"
$ws.Range("E25").Value = "'1.12%"
$ws.Range("D26").Value = "'0.0001293"
$ws.Range("E26").Value = "'-0.67%"
$ws.Range("D27").Value = "'0.0003363"
$ws.Range("E27").Value = "'-1.02%"
$ws.Range("D39").Value = "'0.01717"
$ws.Range("E39").Value = "'7.33%"
$ws.Range("D40").Value = "'0.04726"
$ws.Range("E40").Value = "'6.38%"
$ws.Range("D41").Value = "'0.007742"
$ws.Range("E41").Value = "'5.51%"
$ws.Range("D42").Value = "'0.009742"
$ws.Range("E42").Value = "'11.26%"
$ws.Range("D43").Value = "'0.1404"
$ws.Range("E43").Value = "'5.79%"
$ws.Range("D44").Value = "'0.002065"
$ws.Range("E44").Value = "'8.56%"
$ws.Range("D45").Value = "'0.009682"
$ws.Range("E45").Value = "'5.23%"
$ws.Range("D46").Value = "'0.00006055"
$ws.Range("E46").Value = "'1.56%"
$ws.Range("D47").Value = "'0.00000000743"
$ws.Range("E47").Value = "'-1.06%"
$ws.Range("D48").Value = "'2.551"
$ws.Range("E48").Value = "'13.84%"
$ws.Range("D49").Value = "'0.001982"
$ws.Range("E49").Value = "'-31.64%"
$ws.Range("D50").Value = "'0.00002081"
$ws.Range("E50").Value = "'-1.06%"
$ws.Range("D51").Value = "'0.0001982"
$ws.Range("E51").Value = "'-1.06%"
